# Apply "More tweaks and de-dups" changes to the BINGO workbook.
# Target sheet: "list" (sheet2.xml) - the master phrase list in column A.
#
# Changes:
#   1. Remove the row containing "Dog barking" entirely (rows below shift up).
#   2. Reword "It's on my radar"                 -> "on the radar"
#   3. Lower-case "Wind / road noise"            -> "wind / road noise"
#   4. Lower-case "Typing noise"                 -> "typing noise"
#   5. Lower-case the "Unnecessary verbing (...)" phrase
#
# The text edits are issued in this specific order so new shared-string
# entries get appended to the workbook's string table in the same order
# as the canonical edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

# 1. Delete the "Dog barking" row (row 12) - entire row, shifting cells up.
$ws.Rows.Item(12).Delete()

# 2. "It's on my radar" was row 47, now row 46; reworded to "on the radar".
$ws.Cells.Item(46, 1).Value = "on the radar"

# 3. "Wind / road noise" was row 19, now row 18 after the deletion above.
$ws.Cells.Item(18, 1).Value = "wind / road noise"

# 4. "Typing noise" is row 6; unaffected by the row deletion.
$ws.Cells.Item(6, 1).Value = "typing noise"

# 5. Unnecessary verbing phrase was row 46, now row 45.
$ws.Cells.Item(45, 1).Value = "unnecessary verbing (""Let's solution that"")"

# Match the author's final cursor position/selection in the saved file.
$ws.Activate() | Out-Null
$ws.Range("A52").Select() | Out-Null
